$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.326.76"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.871.59"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "235.17"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.4699"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "0.2872"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "0.06589"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "21.69"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "0.07910"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "1.866.90"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "0.6924"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "5.111"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").Value = "268.82"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "30.297.76"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "0.000007691"
$ws.Range("E19").Value = "  +3.80%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "2.117.60"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D23").Value = "5.257"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").Value = "6.211"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "9.396"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "167.48"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "18.91"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "1.360"
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").Value = "0.09891"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").Value = "4.364"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").Value = "4.073"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "0.04757"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "0.7046"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").Value = "2.804"
$ws.Range("E39").Value = "  +6.41%  "
$ws.Range("D40").Value = "6.230"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Value = "73.05"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("D42").Value = "1.956"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Value = "0.4181"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "102.75"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "7.148"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "942.03"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").Value = "9.174"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("E51").Value = "  +0.51%  "
